$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 22 de Septiembre de 2020 a las 18:10"

# Country (column A) relabels caused by re-sorting the table by updated Casos totales
$ws.Cells.Item(44, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(45, 1).Value = "Guatemala"
$ws.Cells.Item(115, 1).Value = "Jordania"
$ws.Cells.Item(116, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(117, 1).Value = "Suazilandia"
$ws.Cells.Item(118, 1).Value = "Cabo Verde"
$ws.Cells.Item(120, 1).Value = "Cuba"
$ws.Cells.Item(121, 1).Value = "Hong Kong"
$ws.Cells.Item(214, 1).Value = "Montserrat"
$ws.Cells.Item(215, 1).Value = "Islas Malvinas"

# Updated statistic values (columns B-H) from refreshed data
$ws.Cells.Item(4, 2).Value = 7053783
$ws.Cells.Item(4, 3).Value = 7567
$ws.Cells.Item(4, 4).Value = 4301523
$ws.Cells.Item(4, 5).Value = 2547459
$ws.Cells.Item(4, 7).Value = 295
$ws.Cells.Item(4, 8).Value = 204801
$ws.Cells.Item(17, 2).Value = 403551
$ws.Cells.Item(17, 3).Value = 4926
$ws.Cells.Item(17, 7).Value = 37
$ws.Cells.Item(17, 8).Value = 41825
$ws.Cells.Item(23, 2).Value = 300897
$ws.Cells.Item(23, 3).Value = 1392
$ws.Cells.Item(23, 4).Value = 219670
$ws.Cells.Item(23, 5).Value = 45489
$ws.Cells.Item(23, 7).Value = 14
$ws.Cells.Item(23, 8).Value = 35738
$ws.Cells.Item(25, 2).Value = 276504
$ws.Cells.Item(25, 3).Value = 953
$ws.Cells.Item(25, 5).Value = 20717
$ws.Cells.Item(25, 7).Value = 6
$ws.Cells.Item(25, 8).Value = 9487
$ws.Cells.Item(29, 2).Value = 146385
$ws.Cells.Item(29, 3).Value = 970
$ws.Cells.Item(29, 4).Value = 126230
$ws.Cells.Item(29, 5).Value = 10923
$ws.Cells.Item(29, 7).Value = 4
$ws.Cells.Item(29, 8).Value = 9232
$ws.Cells.Item(44, 2).Value = 86447
$ws.Cells.Item(44, 3).Value = 852
$ws.Cells.Item(44, 4).Value = 76025
$ws.Cells.Item(44, 5).Value = 10017
$ws.Cells.Item(44, 8).Value = 405
$ws.Cells.Item(45, 2).Value = 85681
$ws.Cells.Item(45, 4).Value = 75172
$ws.Cells.Item(45, 5).Value = 7385
$ws.Cells.Item(45, 8).Value = 3124
$ws.Cells.Item(55, 5).Value = 6899
$ws.Cells.Item(55, 7).Value = 3
$ws.Cells.Item(55, 8).Value = 227
$ws.Cells.Item(61, 4).Value = 42100
$ws.Cells.Item(61, 5).Value = 6510
$ws.Cells.Item(64, 2).Value = 47446
$ws.Cells.Item(64, 3).Value = 650
$ws.Cells.Item(64, 4).Value = 35542
$ws.Cells.Item(64, 5).Value = 10674
$ws.Cells.Item(64, 7).Value = 19
$ws.Cells.Item(64, 8).Value = 1230
$ws.Cells.Item(70, 2).Value = 37218
$ws.Cells.Item(70, 3).Value = 139
$ws.Cells.Item(70, 4).Value = 24147
$ws.Cells.Item(70, 5).Value = 12412
$ws.Cells.Item(70, 7).Value = 9
$ws.Cells.Item(70, 8).Value = 659
$ws.Cells.Item(88, 2).Value = 15928
$ws.Cells.Item(88, 3).Value = 333
$ws.Cells.Item(88, 5).Value = 5587
$ws.Cells.Item(88, 7).Value = 8
$ws.Cells.Item(88, 8).Value = 352
$ws.Cells.Item(94, 2).Value = 12666
$ws.Cells.Item(94, 3).Value = 131
$ws.Cells.Item(94, 4).Value = 7042
$ws.Cells.Item(94, 5).Value = 5257
$ws.Cells.Item(115, 2).Value = 5679
$ws.Cells.Item(115, 3).Value = 634
$ws.Cells.Item(115, 4).Value = 3707
$ws.Cells.Item(115, 5).Value = 1939
$ws.Cells.Item(115, 7).Value = 1
$ws.Cells.Item(115, 8).Value = 33
$ws.Cells.Item(116, 2).Value = 5404
$ws.Cells.Item(116, 4).Value = 5336
$ws.Cells.Item(116, 5).Value = 7
$ws.Cells.Item(116, 8).Value = 61
$ws.Cells.Item(117, 2).Value = 5282
$ws.Cells.Item(117, 4).Value = 4647
$ws.Cells.Item(117, 5).Value = 531
$ws.Cells.Item(117, 8).Value = 104
$ws.Cells.Item(118, 2).Value = 5281
$ws.Cells.Item(118, 3).Value = 0
$ws.Cells.Item(118, 4).Value = 4674
$ws.Cells.Item(118, 5).Value = 555
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 52
$ws.Cells.Item(119, 2).Value = 5270
$ws.Cells.Item(119, 3).Value = 127
$ws.Cells.Item(119, 4).Value = 1444
$ws.Cells.Item(119, 5).Value = 3751
$ws.Cells.Item(119, 7).Value = 5
$ws.Cells.Item(119, 8).Value = 75
$ws.Cells.Item(120, 2).Value = 5222
$ws.Cells.Item(120, 3).Value = 81
$ws.Cells.Item(120, 4).Value = 4506
$ws.Cells.Item(120, 5).Value = 599
$ws.Cells.Item(120, 7).Value = 1
$ws.Cells.Item(120, 8).Value = 117
$ws.Cells.Item(121, 2).Value = 5047
$ws.Cells.Item(121, 3).Value = 8
$ws.Cells.Item(121, 4).Value = 4717
$ws.Cells.Item(121, 5).Value = 227
$ws.Cells.Item(121, 8).Value = 103
$ws.Cells.Item(130, 2).Value = 3974
$ws.Cells.Item(130, 3).Value = 29
$ws.Cells.Item(130, 4).Value = 1871
$ws.Cells.Item(130, 5).Value = 2038
$ws.Cells.Item(160, 2).Value = 1618
$ws.Cells.Item(160, 3).Value = 15
$ws.Cells.Item(160, 5).Value = 227
$ws.Cells.Item(195, 2).Value = 115
$ws.Cells.Item(195, 3).Value = 1
$ws.Cells.Item(195, 5).Value = 4
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 8).Value = 1
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 8).Value = 0
